$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B:E
$ws.Range("B2").Value = 9.4573355776025316
$ws.Range("C2").Value = 5.5912823976564896
$ws.Range("D2").Value = 9.8323984339473682
$ws.Range("E2").Value = 8.3264446056587449

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 5.7132292175159396
$ws.Range("C3").Value = 7.6542774238505444
$ws.Range("D3").Value = 6.6232210289869027
$ws.Range("E3").Value = 8.4670984845031541

# Reflect the new active selection, matching the edited range
$ws.Range("B1:E3").Select()
